# Add a new "Dinsdag" (Tuesday) entry to the Logboek Joey timesheet (row 33),
# mirroring the existing rows (e.g. row 30) for the "Week 3" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day name
$ws.Range("B33").Value = "Dinsdag"

# Begin / end time, formatted as time-of-day like the other rows
$ws.Range("C33").Value = 0.54166666666666663
$ws.Range("C33").NumberFormat = "h:mm"

$ws.Range("D33").Value = 0.61458333333333337
$ws.Range("D33").NumberFormat = "h:mm"

# Minutes worked
$ws.Range("E33").Value = 330

# Description of the work performed
$ws.Range("G33").Value = "aan c# voor applicate gewerkt"

# Move the on-screen selection/scroll position to just below the new row,
# matching where the author's cursor ended up after the edit.
$ws.Range("G34").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
